$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.312.59"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.588.48"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "1.598.89"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "26.312.83"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "1.314.28"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.611"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -12.11%  "
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "1.724.63"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0981"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.81%  "
